$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 and add the new sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Type of quanv"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Diff num of quanv filter"

# --- Column widths on the new sheet ---
# (ColumnWidth is quantized to 1/7-character pixel steps by the engine; these
# are the closest achievable values to the authored 35.125 / 23.125 widths)
$ws2.Range("A1").ColumnWidth = 34.428571428571428
$ws2.Range("B1:C1").ColumnWidth = 22.428571428571428

# --- Labels ---
$ws2.Range("A2").Value = "MNIST"
$ws2.Range("A4").Value = "Number of quanv filter"

# --- Fix the typo on Sheet1 (G4: "Not nunned" -> "Not runned") ---
# (done after the new-sheet strings so the shared-string table order matches)
$ws1.Range("G4").Value = "Not runned"

# --- Column A values (filter counts) ---
$ws2.Range("A5").Value = 1
$ws2.Range("A6").Value = 2
$ws2.Range("A7").Value = 3
$ws2.Range("A8").Value = 4
$ws2.Range("A9").Value = 5
$ws2.Range("A10").Value = 6
$ws2.Range("A11").Value = 7
$ws2.Range("A12").Value = 8
$ws2.Range("A13").Value = 9
$ws2.Range("A14").Value = 10
$ws2.Range("A15").Value = 20
$ws2.Range("A16").Value = 30
$ws2.Range("A17").Value = 40
$ws2.Range("A18").Value = 50
$ws2.Range("A19").Value = 100

# --- Column B values (status number) ---
$ws2.Range("B5").Value = 2
$ws2.Range("B6").Value = 2
$ws2.Range("B7").Value = 1
$ws2.Range("B8").Value = 1
$ws2.Range("B9").Value = 2
$ws2.Range("B10").Value = 1
$ws2.Range("B11").Value = 1
$ws2.Range("B14").Value = 2
$ws2.Range("B15").Value = 2
$ws2.Range("B16").Value = 2
$ws2.Range("B17").Value = 2
$ws2.Range("B18").Value = 2
$ws2.Range("B19").Value = 2

# --- Legend block (D/E columns) ---
$ws2.Range("D6").Value = "Status"
$ws2.Range("D14").Value = 0
$ws2.Range("E14").Value = "Not runned"
$ws2.Range("D15").Value = 1
$ws2.Range("E15").Value = "Running"
$ws2.Range("D16").Value = 2
$ws2.Range("E16").Value = "Runned"
$ws2.Range("D18").Value = "PC"
$ws2.Range("D19").Value = "haivt"
$ws2.Range("D20").Value = "lantran"
$ws2.Range("D21").Value = "tinhuynh"
$ws2.Range("D22").Value = "HPC"

# --- Formatting: reuse existing fill styles from Sheet1 via copy/paste-special ---
# Note: B5,B9,B12,B13,B14,B16,B17,B18 use a "no-fill" style (fillId=0, just an
# applyFill flag) in the source workbook - visually identical to no style at
# all, so it is intentionally not recreated here.

# s="3" (yellow) -> B6,B7,B8,B15 and D21
$ws1.Range("F11").Copy()
$ws2.Range("B6").PasteSpecial(-4122)
$ws2.Range("B7").PasteSpecial(-4122)
$ws2.Range("B8").PasteSpecial(-4122)
$ws2.Range("B15").PasteSpecial(-4122)
$ws2.Range("D21").PasteSpecial(-4122)

# s="2" (theme accent6 green) -> B10,B11 and D20
$ws1.Range("F10").Copy()
$ws2.Range("B10").PasteSpecial(-4122)
$ws2.Range("B11").PasteSpecial(-4122)
$ws2.Range("D20").PasteSpecial(-4122)

# s="1" (theme accent2 orange) -> B19 and D19
$ws1.Range("F9").Copy()
$ws2.Range("B19").PasteSpecial(-4122)
$ws2.Range("D19").PasteSpecial(-4122)

# s="4" (red) -> D22
$ws1.Range("F12").Copy()
$ws2.Range("D22").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Selections ---
$ws1.Range("I7").Select()
$ws2.Range("D16").Select()
